$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# New headers for the "usage" feature (replacing MinProduction/MaxProduction
# with IsSet/UsageMin, and adding UsageAvg/UsageMax columns)
$ws.Range("H1").Value2 = "IsSet"
$ws.Range("I1").Value2 = "UsageMin"
$ws.Range("J1").Value2 = "UsageAvg"
$ws.Range("K1").Value2 = "UsageMax"

# Row 2 (CrOil)
$ws.Range("H2").Value2 = 1
$ws.Range("I2").Value2 = 200
$ws.Range("J2").Value2 = 300
$ws.Range("K2").Value2 = 400

# Row 3 (H2)
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0

# Row 4 (BM)
$ws.Range("I4").Value2 = 10
$ws.Range("J4").Value2 = 660000
$ws.Range("K4").Value2 = 500

# Sources becomes the active sheet / tab again, with J4 selected
$ws.Select()
$ws.Range("J4").Select()
